$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily data rows for 2021-05-28 through 2021-06-28
# Columns: A = date serial, B = nuovi pos., C = somma mobile 7gg., D = somma mobile 7gg. per 100mila abitanti
$data = @(
    @(44344,1,1,12.15362177929023),
    @(44345,0,1,12.15362177929023),
    @(44346,2,3,36.46086533787069),
    @(44347,1,4,48.61448711716091),
    @(44348,1,5,60.76810889645115),
    @(44349,0,5,60.76810889645115),
    @(44350,2,7,85.07535245503161),
    @(44351,0,6,72.92173067574137),
    @(44352,0,6,72.92173067574137),
    @(44353,3,7,85.07535245503161),
    @(44354,1,7,85.07535245503161),
    @(44355,2,8,97.22897423432183),
    @(44356,0,8,97.22897423432183),
    @(44357,1,7,85.07535245503161),
    @(44358,1,8,97.22897423432183),
    @(44359,1,9,109.3825960136121),
    @(44360,0,6,72.92173067574137),
    @(44361,0,5,60.76810889645115),
    @(44362,0,3,36.46086533787069),
    @(44363,0,3,36.46086533787069),
    @(44364,1,3,36.46086533787069),
    @(44365,0,2,24.30724355858046),
    @(44366,0,1,12.15362177929023),
    @(44367,0,1,12.15362177929023),
    @(44368,0,1,12.15362177929023),
    @(44369,0,1,12.15362177929023),
    @(44370,0,1,12.15362177929023),
    @(44371,0,0,0),
    @(44372,0,0,0),
    @(44373,0,0,0),
    @(44374,0,0,0),
    @(44375,0,0,0)
)

$startRow = 270

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}

$endRow = $startRow + $data.Count - 1

# Match the date formatting/style of column A (copied from the last pre-existing
# data row, A269) onto all the newly added A270:A301 cells
$ws.Range("A269").Copy()
$ws.Range("A$startRow`:A$endRow").PasteSpecial(-4122)
$excel.CutCopyMode = $false
